$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume 1h %) updates per latest scrape.
# A leading "'" forces text interpretation for values that would
# otherwise be auto-parsed as numbers (e.g. "683.12"), matching how
# Excel treats a user keystroke of '683.12 into a General cell.

$ws.Range("D2").Value = "69.299.84"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "3.691.37"
$ws.Range("E3").Value = "  -2.87%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'683.12"
$ws.Range("E5").Value = "  -2.98%  "
$ws.Range("E6").Value = "  -4.53%  "
$ws.Range("D7").Value = "3.687.57"
$ws.Range("E7").Value = "  -2.99%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  -3.91%  "
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "  -7.01%  "
$ws.Range("D11").Value = "'7.28"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("E13").Value = "  -6.29%  "
$ws.Range("D14").Value = "'33.56"
$ws.Range("E14").Value = "  -7.04%  "
$ws.Range("D15").Value = "4.316.09"
$ws.Range("E15").Value = "  -2.84%  "
$ws.Range("D16").Value = "3.701.74"
$ws.Range("E16").Value = "  -4.06%  "
$ws.Range("D17").Value = "69.413.92"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").Value = "'16.35"
$ws.Range("E19").Value = "  -5.88%  "
$ws.Range("E20").Value = "  -7.24%  "
$ws.Range("D21").Value = "'484.50"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "'9.79"
$ws.Range("E22").Value = "  -7.75%  "
$ws.Range("D23").Value = "'0.667"
$ws.Range("E23").Value = "  -8.28%  "
$ws.Range("D24").Value = "'80.05"
$ws.Range("E24").Value = "  -5.64%  "
$ws.Range("D25").Value = "3.837.21"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("D26").Value = "'0.0000128"
$ws.Range("E26").Value = "  -10.87%  "
$ws.Range("D27").Value = "'11.54"
$ws.Range("E27").Value = "  -4.63%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'9.59"
$ws.Range("E29").Value = "  -8.07%  "
$ws.Range("D30").Value = "'1.84"
$ws.Range("E30").Value = "  -10.13%  "
$ws.Range("D31").Value = "'2.75"
$ws.Range("E31").Value = "  -10.58%  "
$ws.Range("D32").Value = "'2.12"
$ws.Range("E32").Value = "  -4.17%  "
$ws.Range("D33").Value = "'6.82"
$ws.Range("E33").Value = "  -6.73%  "
$ws.Range("D34").Value = "'27.10"
$ws.Range("E34").Value = "  -6.70%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  -4.40%  "
$ws.Range("D37").Value = "3.655.75"
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("D38").Value = "'8.57"
$ws.Range("E38").Value = "  -5.63%  "
$ws.Range("D39").Value = "'6.03"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "'0.0942"
$ws.Range("E40").Value = "  -6.91%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -5.99%  "
$ws.Range("E43").Value = "  -7.16%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'157.67"
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D48").Value = "'0.000279"
$ws.Range("E48").Value = "  -12.99%  "
$ws.Range("D49").Value = "'390.84"
$ws.Range("E49").Value = "  -8.23%  "
$ws.Range("D50").Value = "'8.10"
$ws.Range("E50").Value = "  -5.77%  "
$ws.Range("E51").Value = "  -5.10%  "
